# Updated symbol list on Sun Jan 15 08:58:07 UTC 2023 with GitHub Actions
# Applies the coinranking.com data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value.
# All values in columns D/E are numeric-looking strings (prices / percentages)
# that must stay stored as literal TEXT (matching the original inlineStr cells),
# so we write them with a leading apostrophe to suppress Excel's auto number/
# percentage conversion, then reset the cell style back to "Normal" so no stray
# number-format style is left behind (the sheet keeps its original formatting).
$updates = @(
    @{ Cell = 'D2'; Value = '295.22' }
    @{ Cell = 'E2'; Value = '-5.75%' }
    @{ Cell = 'D3'; Value = '31.44' }
    @{ Cell = 'E3'; Value = '-3.56%' }
    @{ Cell = 'E4'; Value = '-4.22%' }
    @{ Cell = 'D5'; Value = '0.07443' }
    @{ Cell = 'D6'; Value = '7.717' }
    @{ Cell = 'E6'; Value = '-2.26%' }
    @{ Cell = 'B7'; Value = 'FTXToken' }
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' }
    @{ Cell = 'D7'; Value = '1.697' }
    @{ Cell = 'E7'; Value = '4.91%' }
    @{ Cell = 'B8'; Value = 'GateToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D8'; Value = '3.781' }
    @{ Cell = 'E8'; Value = '0.86%' }
    @{ Cell = 'D9'; Value = '0.9296' }
    @{ Cell = 'D10'; Value = '0.1686' }
    @{ Cell = 'E10'; Value = '-2.28%' }
    @{ Cell = 'D11'; Value = '0.07135' }
    @{ Cell = 'E11'; Value = '-5.96%' }
    @{ Cell = 'D12'; Value = '0.07934' }
    @{ Cell = 'E12'; Value = '-4.53%' }
    @{ Cell = 'D13'; Value = '0.02998' }
    @{ Cell = 'E13'; Value = '-1.25%' }
    @{ Cell = 'D14'; Value = '0.09916' }
    @{ Cell = 'E14'; Value = '0.06%' }
    @{ Cell = 'D15'; Value = '0.001489' }
    @{ Cell = 'E15'; Value = '-2.49%' }
    @{ Cell = 'D16'; Value = '0.006216' }
    @{ Cell = 'E16'; Value = '0.51%' }
    @{ Cell = 'D18'; Value = '2.225' }
    @{ Cell = 'E18'; Value = '-0.92%' }
    @{ Cell = 'E19'; Value = '-1.31%' }
    @{ Cell = 'E20'; Value = '0.10%' }
    @{ Cell = 'D21'; Value = '4.560' }
    @{ Cell = 'E21'; Value = '7.54%' }
    @{ Cell = 'D22'; Value = '0.04650' }
    @{ Cell = 'E22'; Value = '1.76%' }
    @{ Cell = 'D24'; Value = '0.001219' }
    @{ Cell = 'E24'; Value = '-0.24%' }
    @{ Cell = 'D25'; Value = '0.004416' }
    @{ Cell = 'E25'; Value = '-2.05%' }
    @{ Cell = 'D26'; Value = '0.0001303' }
    @{ Cell = 'E26'; Value = '0.38%' }
    @{ Cell = 'E27'; Value = '8.03%' }
    @{ Cell = 'D39'; Value = '0.01661' }
    @{ Cell = 'E39'; Value = '-6.05%' }
    @{ Cell = 'D40'; Value = '0.04431' }
    @{ Cell = 'E40'; Value = '-4.68%' }
    @{ Cell = 'D41'; Value = '0.007076' }
    @{ Cell = 'E41'; Value = '-1.62%' }
    @{ Cell = 'D42'; Value = '0.1326' }
    @{ Cell = 'E42'; Value = '-3.49%' }
    @{ Cell = 'D43'; Value = '0.002093' }
    @{ Cell = 'E43'; Value = '-7.29%' }
    @{ Cell = 'D44'; Value = '0.01231' }
    @{ Cell = 'E44'; Value = '-14.70%' }
    @{ Cell = 'D45'; Value = '0.00006013' }
    @{ Cell = 'E45'; Value = '-3.18%' }
    @{ Cell = 'D46'; Value = '0.7116' }
    @{ Cell = 'E46'; Value = '-62.40%' }
    @{ Cell = 'D47'; Value = '0.01102' }
    @{ Cell = 'E47'; Value = '-15.15%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $text = $u.Value
    if ($text -match '^-?[0-9.]+%?$') {
        # Numeric-looking value (price or percentage) -> force text.
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    } else {
        # Plain text (coin name / URL) -> assign directly.
        $cell.Value = $text
    }
}
